# "Horas de Produção" — fill in the missing week (row 7, date 43738 = 09/30/2019)
# with its worked-hours entries, and add a weekly-totals row (row 18) that sums
# each person's column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 currently holds only the date (A7); give B7:F7 the same number format /
# borders as the other weekly rows (copy formatting from row 4) before writing
# the actual time values, then match G7's "red note" style from the row above.
$ws.Range("B4:F4").Copy()
$ws.Range("B7:F7").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("G6").Copy()
$ws.Range("G7").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Range("B7").Value = 0.15625
$ws.Range("C7").Value = 0.15347222222222223
$ws.Range("D7").Value = 0.15486111111111112
$ws.Range("E7").Value = 0.14791666666666667
$ws.Range("F7").Value = 0.15625
$ws.Range("G7").Value = "-20 min de intervalo"

# Row 18: totals per column across the whole table (B4:B17 … F4:F17), shown as hh:mm.
$ws.Range("B18:F18").Formula = "=sum(B4:B17)"
$ws.Range("B18:F18").NumberFormat = "hh:mm"
